$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.467.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6317"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07564"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2961"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.65"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07708"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.83"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6836"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.148"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.497.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.90"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.528"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.78"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1399"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.386"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.468"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.275"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05686"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.129"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.846"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.159"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7156"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.587"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.247.10"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01808"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.780"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9109"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.212"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.69"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000121"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.071"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.136"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4027"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.694"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1122"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05721"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.73%  "
